# Aashish Sharma - Task Breakdown sheet update
# Story SSDMS-51 (rows 21-35): several tasks moved from "In-Dev"/"Not-Started" to
# "Done" - their "Hours Burnt" (F) are bumped up to match "Planned Hours" (E),
# so "Remaining Hours" (G, = E-F) drops to 0, and the Status (H) cell now reads
# "Done" (picking up the fill/format already used by the other "Done" rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "Hours Burnt" (column F) for the tasks that got finished ----
$ws.Range("F27").Value = 4
$ws.Range("F28").Value = 6
$ws.Range("F30").Value = 6
$ws.Range("F32").Value = 2
$ws.Range("F33").Value = 1
$ws.Range("F34").Value = 1
$ws.Range("F35").Value = 2
# (column G holds =E-F formulas and recalculates automatically)

# --- 2. Flip the Status (column H) cells to "Done" ---------------------
$doneRows = @(27, 28, 29, 30, 32, 33, 34, 35, 36)
foreach ($r in $doneRows) {
    $ws.Range("H$r").Value = "Done"
}

# Re-use the fill/format already applied to the existing "Done" status cells
# (e.g. H21) so the changed cells pick up the same green shading instead of
# keeping their old "In-Dev"/"Not-Started" fill.
$doneFormat = $ws.Range("H21")
$doneFormat.Copy()
foreach ($r in $doneRows) {
    $ws.Range("H$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- 3. Scroll position / active selection, matching the author's view ----
$ws.Range("H36").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 30
$win.ScrollColumn = 3
